$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" (E) / "Valor Mora" (F) block for rows 16-22 is reversed in
# order: row 16 swaps with row 22, row 17 swaps with row 21, row 18 swaps with
# row 20; row 19 (the middle row) stays where it is.
$topRow = 16
$bottomRow = 22

while ($topRow -lt $bottomRow) {
    $eTop = $ws.Cells.Item($topRow, 5).Value2
    $fTop = $ws.Cells.Item($topRow, 6).Value2

    $eBottom = $ws.Cells.Item($bottomRow, 5).Value2
    $fBottom = $ws.Cells.Item($bottomRow, 6).Value2

    $ws.Cells.Item($topRow, 5).Value2 = $eBottom
    $ws.Cells.Item($topRow, 6).Value2 = $fBottom

    $ws.Cells.Item($bottomRow, 5).Value2 = $eTop
    $ws.Cells.Item($bottomRow, 6).Value2 = $fTop

    $topRow = $topRow + 1
    $bottomRow = $bottomRow - 1
}
